# Applies the "history -> chemistry" content rewrite described by the diff.
$d = $word.ActiveDocument

# NOTE on Find.Execute:
#  * Replace:=wdReplaceAll collapses the Range to the last replaced match, so
#    a Range variable must never be reused for more than one Execute call --
#    always operate on a freshly duplicated Range.
#  * Passing the replacement text straight to Execute() runs it through
#    Word's autocorrect/"smart quotes" pass (straight ' becomes a curly
#    U+2019). To keep straight apostrophes straight, Execute() is called in
#    "find only" mode (no replacement argument) and the match text is then
#    overwritten via Range.Text, which does not autocorrect.

function Replace-Text($paragraphRange, $old, $new) {
    $rng = $paragraphRange.Duplicate
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $rng.Text = $new
    } else {
        throw "Find text not located: $old"
    }
}

# --- Title ---
Replace-Text $d.Paragraphs.Item(1).Range "A Transformative Journey Through Time" "The Symphony of Matter: An Exploration into the Realm of Chemistry"

# --- Author name ---
Replace-Text $d.Paragraphs.Item(2).Range "Marianne Geary" "Emily Williams"

# --- Email (username / domain), scoped to the email paragraph to be safe ---
Replace-Text $d.Paragraphs.Item(3).Range "marianne" "emily"
Replace-Text $d.Paragraphs.Item(3).Range "geary@mail" "williams@schoolmail"
Replace-Text $d.Paragraphs.Item(3).Range "com" "edu"

# --- Body paragraph 1 (paragraph 5) ---

Replace-Text $d.Paragraphs.Item(5).Range "Humanity's relentless pursuit of understanding our place in the grand tapestry of time has fueled groundbreaking advancements in the study of history" "In the vast tapestry of the natural world, chemistry stands as a symphony of matter, beckoning us to unravel its intricacies"

Replace-Text $d.Paragraphs.Item(5).Range " From deciphering ancient scripts and excavating lost civilizations to unraveling the intricacies of human behavior, historians embark on a transformative journey through time, seeking to illuminate the past and its profound impact on our present" " This captivating science delves into the very essence of substances, their composition, and the transformation they undergo"

Replace-Text $d.Paragraphs.Item(5).Range " In this captivating odyssey, scholars meticulously piece together fragments of bygone eras, uncovering forgotten narratives and shedding light on the choices that have shaped our world" " From the smallest subatomic particles to the colossal molecular structures, chemistry unveils the profound connections that govern the universe around us"

# Remove the next two sentences entirely (". Their quest..." + ". Through meticulous analysis...") -
# keeping the trailing period that closes out that run of sentences.
Replace-Text $d.Paragraphs.Item(5).Range ". Their quest for knowledge extends beyond textbooks and libraries, as they venture into archives, museums, and archaeological sites, unearthing tangible remnants of the past. Through meticulous analysis and interpretation, historians weave these fragments into a rich tapestry of understanding, revealing the currents of human civilization that have flowed through the ages" ""

Replace-Text $d.Paragraphs.Item(5).Range "Guided by a deep reverence for the past, historians embark on a quest to unveil the lives of individuals and communities, shedding light on their struggles, triumphs, and enduring legacies" "Chemistry weaves a tale of elements, each possessing a unique identity and story"

Replace-Text $d.Paragraphs.Item(5).Range " They delve into the motivations, beliefs, and values that have shaped the course of history, exploring how these factors have influenced the rise and fall of civilizations, the evolution of societies, and the interplay between different cultures" " These fundamental building blocks combine in myriad ways, forming compounds of infinite variety"

Replace-Text $d.Paragraphs.Item(5).Range " Through their meticulous research and insightful analysis, historians illuminate the intricate web of connections that link past events to our present circumstances, providing invaluable perspectives on the human experience" " As atoms dance and electrons waltz, new substances emerge, exhibiting properties that astound and inspire"

Replace-Text $d.Paragraphs.Item(5).Range " By peering into the mirror of time, we gain a profound appreciation for the complexities of our shared history and the enduring impact it continues to have on our lives" " Chemistry grants us the power to manipulate and harness these elements, synthesizing materials with bespoke properties, from gleaming metals to life-saving drugs"

Replace-Text $d.Paragraphs.Item(5).Range "Moreover, the study of history cultivates critical thinking skills, enabling us to evaluate information, identify biases, and form informed judgments" "The study of chemistry is not merely an intellectual pursuit; it is an art form, a symphony of discovery and creation"

Replace-Text $d.Paragraphs.Item(5).Range " By examining multiple perspectives and engaging in thoughtful discourse, we develop the ability to navigate the complexities of the world around us, making more informed decisions and fostering a deeper understanding of the challenges and opportunities we face" " Chemists wield their knowledge as artists wield their brushes, transforming raw materials into masterpieces of molecular architecture"

# Final sentence gets rewritten and two new sentences are spliced in before the trailing period.
Replace-Text $d.Paragraphs.Item(5).Range " History's lessons serve as a valuable guide, reminding us of the pitfalls and triumphs of past generations, empowering us to forge a better future for ourselves and succeeding generations." " They unlock the secrets of nature, deciphering the intricate language of chemical reactions. Through experimentation and innovation, they orchestrate new compounds, choreographing the dance of atoms to produce substances that benefit humanity."

# --- "Summary" heading paragraph: touch it so the stale lastRenderedPageBreak cache is dropped ---
Replace-Text $d.Paragraphs.Item(6).Range "Summary" "Summary"

# --- Summary body paragraph (paragraph 7) ---

Replace-Text $d.Paragraphs.Item(7).Range "Through the study of history, we embark on a transformative journey through time, unearthing forgotten narratives and illuminating the profound impact of the past on our present" "Chemistry, a symphony of matter, unveils the universe's profound connections"

Replace-Text $d.Paragraphs.Item(7).Range " Historians delve into the lives of individuals and communities, deciphering ancient scripts and excavating lost civilizations to unravel the intricacies of human civilization" " It weaves a tale of elements, their combinations, and transformations, empowering us to manipulate and harness them"

Replace-Text $d.Paragraphs.Item(7).Range " Their meticulous research and insightful analysis provide invaluable perspectives on the human experience, fostering critical thinking skills and enabling us to make informed decisions" " Chemists, as artists, wield their knowledge to orchestrate new compounds, benefiting humanity"

Replace-Text $d.Paragraphs.Item(7).Range " History's lessons guide us, reminding us of past triumphs and pitfalls, empowering us to forge a better future" " The study of chemistry is an art form, deciphering nature's language and choreographing the dance of atoms to create substances that improve our world"

# --- Add a new trailing empty paragraph at the end of the document (before sectPr) ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null
